# Apply attendance-count updates to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> list of column letters that should become 1
$updates = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("H")
    6  = @("D", "E")
    7  = @("H")
    8  = @("H")
    9  = @("H")
    10 = @("G", "H")
    11 = @("D", "E")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("H")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
